$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "VALOR MORA" total (E11) and "Cant. Periodos" count (F13) ---
$ws.Range("E11").Value = 341628
$ws.Range("F13").Value = 3

# --- 2. Center the "Periodo Mora" column (E) for the existing data rows ---
$ws.Range("E16:E19").HorizontalAlignment = -4108

# --- 3. Insert two new rows for the new "2509" period, right after the existing
#        2508 block (rows 18-19), pushing the signature-block rows down. ---
$ws.Rows("20:21").Insert()

# Duplicate the last existing data block (rows 18:19, the "2508" pair) down into
# the newly inserted rows 20:21, carrying values + formatting (borders etc.)
# so the new block ends up with the same "closing" border style that used to
# belong to row 19.
$ws.Range("B18:J19").Copy($ws.Range("B20:J21"))

# Row 19 is no longer the last row of the table, so it should pick up the
# "interior" row formatting (like row 18) instead of the "closing" border
# formatting it had before - copy formats only (xlPasteFormats), values stay.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# --- 4. Fill in the new "2509" period row values (Periodo Mora column E) ---
$ws.Range("E20").Value = "2509"
$ws.Range("E21").Value = "2509"
